# Add the 2022-Q4 data: a new row in the "总计" (totals) summary sheet, and a
# new "2022-Q4" worksheet (fund holdings detail) placed right after "总计" and
# before "2022-Q3". All existing quarterly sheets shift right by one position
# but keep their own data unchanged.

$wb = $excel.ActiveWorkbook

function Set-TextCell($range, $value) {
    # Force the cell to be stored as text (matches the source data, where
    # numeric-looking values such as "1.06" or fund codes like "012751" are
    # kept as literal strings rather than numbers), while leaving the cell's
    # effective formatting/style untouched.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ----------------------------------------------------------------------
# 1) "总计" summary sheet: insert the 2022-Q4 row at the top of the data
#    (row 2), pushing every other quarter down by one row.
# ----------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 7 is brand new (the sheet used to end at row 6) — clone row 6's
# formatting into it first so the new row picks up the same styling
# (bordered/centered index cell in column A) as every other data row.
$summary.Range("A6:D6").Copy($summary.Range("A7:D7")) | Out-Null

$summary.Range("A2").Value = 0
Set-TextCell $summary.Range("B2") "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.09

$summary.Range("A3").Value = 1
Set-TextCell $summary.Range("B3") "2022-Q3"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.1

$summary.Range("A4").Value = 2
Set-TextCell $summary.Range("B4") "2022-Q2"
$summary.Range("C4").Value = 4
$summary.Range("D4").Value = 0.06

$summary.Range("A5").Value = 3
Set-TextCell $summary.Range("B5") "2022-Q1"
$summary.Range("C5").Value = 4
$summary.Range("D5").Value = 0.07

$summary.Range("A6").Value = 4
Set-TextCell $summary.Range("B6") "2021-Q4"
$summary.Range("C6").Value = 3
$summary.Range("D6").Value = 1.2

$summary.Range("A7").Value = 5
Set-TextCell $summary.Range("B7") "2021-Q3"
$summary.Range("C7").Value = 4
$summary.Range("D7").Value = 0.02

# ----------------------------------------------------------------------
# 2) Create the new "2022-Q4" fund-detail sheet. Copying an existing
#    quarterly sheet (same 7-row / 8-column layout and styling) as a
#    template is the simplest way to reproduce its formatting exactly,
#    then the copy's values are overwritten below.
# ----------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$afterAnchor = $wb.Worksheets.Item("总计")
$template.Copy($null, $afterAnchor) | Out-Null
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

Set-TextCell $newSheet.Range("B2") "012751"
Set-TextCell $newSheet.Range("C2") "建信纳斯达克100指数（QDII）A 美元现汇"
Set-TextCell $newSheet.Range("D2") "1.06"
Set-TextCell $newSheet.Range("E2") "82.28"
Set-TextCell $newSheet.Range("F2") "2.54"
Set-TextCell $newSheet.Range("G2") "0.0269"
$newSheet.Range("H2").Value = 7

Set-TextCell $newSheet.Range("B3") "012752"
Set-TextCell $newSheet.Range("C3") "建信纳斯达克100指数（QDII）C 人民币"
Set-TextCell $newSheet.Range("D3") "1.06"
Set-TextCell $newSheet.Range("E3") "82.28"
Set-TextCell $newSheet.Range("F3") "2.54"
Set-TextCell $newSheet.Range("G3") "0.0269"
$newSheet.Range("H3").Value = 7

Set-TextCell $newSheet.Range("B4") "012753"
Set-TextCell $newSheet.Range("C4") "建信纳斯达克100指数（QDII）C 美元现汇"
Set-TextCell $newSheet.Range("D4") "1.06"
Set-TextCell $newSheet.Range("E4") "82.28"
Set-TextCell $newSheet.Range("F4") "2.54"
Set-TextCell $newSheet.Range("G4") "0.0269"
$newSheet.Range("H4").Value = 7

Set-TextCell $newSheet.Range("B5") "539002"
Set-TextCell $newSheet.Range("C5") "建信新兴市场优选混合（QDII）"
Set-TextCell $newSheet.Range("D5") "0.21"
Set-TextCell $newSheet.Range("E5") "73.13"
Set-TextCell $newSheet.Range("F5") "2.76"
Set-TextCell $newSheet.Range("G5") "0.0058"
$newSheet.Range("H5").Value = 8
